$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's win totals as a new row at the bottom of the table
$ws.Range("A58").Value = 46007
$ws.Range("B58").Value = 122
$ws.Range("C58").Value = 137
$ws.Range("D58").Value = 128

# Match the date formatting used by the rest of column A (row 2 through 57)
$ws.Range("A58").NumberFormat = "YYYY-MM-DD HH:MM:SS"
